# Weekly update: insert a new weekly price record (Sweet Heart) before the
# existing "Rainier" row, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 5; rows 5-11 shift down to 6-12.
$ws.Rows.Item(5).EntireRow.Insert()

# Fill in the new row 5 with this week's data.
$ws.Cells.Item(5, 1).Value  = 1
$ws.Cells.Item(5, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(5, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(5, 4).Value  = 44580
$ws.Cells.Item(5, 5).Value  = 15
$ws.Cells.Item(5, 6).Value  = "Fruta"
$ws.Cells.Item(5, 7).Value  = 100103
$ws.Cells.Item(5, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(5, 9).Value  = 100103001
$ws.Cells.Item(5, 10).Value = "Cereza"
$ws.Cells.Item(5, 11).Value = "Sweet Heart"
$ws.Cells.Item(5, 12).Value = "Segunda"
$ws.Cells.Item(5, 13).Value = 300
$ws.Cells.Item(5, 14).Value = 7000
$ws.Cells.Item(5, 15).Value = 8000
$ws.Cells.Item(5, 16).Value = 7500
$ws.Cells.Item(5, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(5, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(5, 19).Value = 750
$ws.Cells.Item(5, 20).Value = 10
